# Auto-generated edit script
# Applies updated market-price figures (columns H-N) for specific leve rows
# across multiple job sheets, per the authoritative diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 7672.778
$ws.Range("J58").Value = 9650.714
$ws.Range("L58").Value = 28952.142
$ws.Range("N58").Value = -29252.142

$ws.Range("H107").Value = 975.6667
$ws.Range("I107").Value = 555
$ws.Range("J107").Value = 1456.4286
$ws.Range("K107").Value = 555
$ws.Range("L107").Value = 1456.4286
$ws.Range("M107").Value = 1365
$ws.Range("N107").Value = -5296.4286

$ws.Range("H113").Value = 1966.6666
$ws.Range("I113").Value = 1975
$ws.Range("K113").Value = 1975
$ws.Range("M113").Value = 1279

$ws.Range("H116").Value = 1966.6666
$ws.Range("I116").Value = 1920
$ws.Range("K116").Value = 1920
$ws.Range("M116").Value = 1522

$ws.Range("H132").Value = 4330829.5
$ws.Range("I132").Value = 4610093
$ws.Range("K132").Value = 13830279
$ws.Range("M132").Value = -13827749

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2454.0264
$ws.Range("I32").Value = 2390.4932
$ws.Range("K32").Value = 2390.4932
$ws.Range("M32").Value = -2103.4932

$ws.Range("H132").Value = 5668.9165
$ws.Range("I132").Value = 9636.25
$ws.Range("J132").Value = 1701.5834
$ws.Range("K132").Value = 28908.75
$ws.Range("L132").Value = 5104.7502
$ws.Range("M132").Value = -26378.75
$ws.Range("N132").Value = -10164.7502

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1881.375
$ws.Range("I105").Value = 846
$ws.Range("K105").Value = 846
$ws.Range("M105").Value = 901

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8697312
$ws.Range("I31").Value = 1722.4286
$ws.Range("K31").Value = 1722.4286
$ws.Range("M31").Value = -1427.4286

$ws.Range("H34").Value = 8697312
$ws.Range("I34").Value = 1722.4286
$ws.Range("K34").Value = 1722.4286
$ws.Range("M34").Value = -1520.4286

$ws.Range("H58").Value = 905.5172
$ws.Range("I58").Value = 849.7895
$ws.Range("J58").Value = 1011.4
$ws.Range("K58").Value = 849.7895
$ws.Range("L58").Value = 1011.4
$ws.Range("M58").Value = -646.7895
$ws.Range("N58").Value = -1417.4

$ws.Range("H99").Value = 1708.8823
$ws.Range("I99").Value = 1400.091
$ws.Range("J99").Value = 2275
$ws.Range("K99").Value = 1400.091
$ws.Range("L99").Value = 2275
$ws.Range("M99").Value = 97.90900000000011
$ws.Range("N99").Value = -5271

$ws.Range("H107").Value = 555.6111
$ws.Range("I107").Value = 514.8095
$ws.Range("J107").Value = 612.73334
$ws.Range("K107").Value = 514.8095
$ws.Range("L107").Value = 612.73334
$ws.Range("M107").Value = 1405.1905
$ws.Range("N107").Value = -4452.73334

$ws.Range("H126").Value = 1708.8823
$ws.Range("I126").Value = 1400.091
$ws.Range("J126").Value = 2275
$ws.Range("K126").Value = 4200.272999999999
$ws.Range("L126").Value = 6825
$ws.Range("M126").Value = -1730.272999999999
$ws.Range("N126").Value = -11765

$ws.Range("H132").Value = 2397.0967
$ws.Range("I132").Value = 1884.95
$ws.Range("J132").Value = 3328.2727
$ws.Range("K132").Value = 5654.85
$ws.Range("L132").Value = 9984.8181
$ws.Range("M132").Value = -3124.85
$ws.Range("N132").Value = -15044.8181

$ws.Range("H134").Value = 1230.5172
$ws.Range("I134").Value = 1118.1875
$ws.Range("J134").Value = 1368.7693
$ws.Range("K134").Value = 3354.5625
$ws.Range("L134").Value = 4106.3079
$ws.Range("M134").Value = -819.5625
$ws.Range("N134").Value = -9176.3079

$ws.Range("H136").Value = 905.5172
$ws.Range("I136").Value = 849.7895
$ws.Range("J136").Value = 1011.4
$ws.Range("K136").Value = 2549.3685
$ws.Range("L136").Value = 3034.2
$ws.Range("M136").Value = 0.63149999999996
$ws.Range("N136").Value = -8134.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = $null
$ws.Range("N42").Value = $null

$ws.Range("H131").Value = 3175638.5
$ws.Range("I131").Value = 1628.8889
$ws.Range("J131").Value = 4274334
$ws.Range("K131").Value = 4886.6667
$ws.Range("L131").Value = 12823002
$ws.Range("M131").Value = 153.3333000000002
$ws.Range("N131").Value = -12833082

$ws.Range("H132").Value = 725.25
$ws.Range("I132").Value = 670.8
$ws.Range("J132").Value = 997.5
$ws.Range("K132").Value = 6037.2
$ws.Range("L132").Value = 8977.5
$ws.Range("M132").Value = -3507.2
$ws.Range("N132").Value = -14037.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16456272
$ws.Range("I70").Value = 18893552
$ws.Range("J70").Value = 4639.5
$ws.Range("K70").Value = 18893552
$ws.Range("L70").Value = 4639.5
$ws.Range("M70").Value = -18893282
$ws.Range("N70").Value = -5179.5

$ws.Range("H73").Value = 16456272
$ws.Range("I73").Value = 18893552
$ws.Range("J73").Value = 4639.5
$ws.Range("K73").Value = 18893552
$ws.Range("L73").Value = 4639.5
$ws.Range("M73").Value = -18892616
$ws.Range("N73").Value = -6511.5

$ws.Range("H80").Value = 4991.3335
$ws.Range("I80").Value = 2198.5715
$ws.Range("J80").Value = 8901.200000000001
$ws.Range("K80").Value = 2198.5715
$ws.Range("L80").Value = 8901.200000000001
$ws.Range("M80").Value = -1200.5715
$ws.Range("N80").Value = -10897.2

$ws.Range("H83").Value = 4991.3335
$ws.Range("I83").Value = 2198.5715
$ws.Range("J83").Value = 8901.200000000001
$ws.Range("K83").Value = 10992.8575
$ws.Range("L83").Value = 44506
$ws.Range("M83").Value = -6000.8575
$ws.Range("N83").Value = -54490

$ws.Range("H132").Value = 54928.895
$ws.Range("I132").Value = 85008.414
$ws.Range("J132").Value = 3364
$ws.Range("K132").Value = 255025.242
$ws.Range("L132").Value = 10092
$ws.Range("M132").Value = -252495.242
$ws.Range("N132").Value = -15152

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1754.1786
$ws.Range("I46").Value = 1507.2858
$ws.Range("J46").Value = 1836.4762
$ws.Range("K46").Value = 1507.2858
$ws.Range("L46").Value = 1836.4762
$ws.Range("M46").Value = -1319.2858
$ws.Range("N46").Value = -2212.4762

$ws.Range("H63").Value = 49271.25
$ws.Range("J63").Value = 49271.25
$ws.Range("L63").Value = 49271.25
$ws.Range("N63").Value = -50769.25

$ws.Range("H66").Value = 49271.25
$ws.Range("J66").Value = 49271.25
$ws.Range("L66").Value = 147813.75
$ws.Range("N66").Value = -155301.75

$ws.Range("H132").Value = 9852.261
$ws.Range("J132").Value = 2999.8572
$ws.Range("L132").Value = 8999.571599999999
$ws.Range("N132").Value = -14059.5716

$ws.Range("H133").Value = 23008.666
$ws.Range("J133").Value = 23008.666
$ws.Range("L133").Value = 23008.666
$ws.Range("N133").Value = -28068.666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 440
$ws.Range("I100").Value = 440
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 880
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -339
$ws.Range("N100").Value = $null

$ws.Range("H126").Value = 25751
$ws.Range("I126").Value = 34001.332
$ws.Range("K126").Value = 102003.996
$ws.Range("M126").Value = -99533.99600000001
